$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header shared-string edits (rich-text run substrings):
#   Volume .. Number 35 -> 36
#   Report Covering the Week 8/26/2024 Through 9/1/2024
#                      -> 9/2/2024 .. 9/8/2024
# ------------------------------------------------------------------
$volRange = $ws.Range("A8")
$volChars = $volRange.Characters(21, 2)
$volChars.Text = "36"

$weekRange = $ws.Range("C9")
$startChars = $weekRange.Characters(27, 9)
$startChars.Text = "9/2/2024"
$endChars = $weekRange.Characters(46, 8)
$endChars.Text = "9/8/2024"

# ------------------------------------------------------------------
# Cells that flip from shared-string text to a genuine number:
# copy the number style/format from a same-style neighbour first,
# then overwrite with the new numeric value (Copy brings over the
# neighbour value+style; the Value= that follows fixes the content).
# ------------------------------------------------------------------
$ws.Range("C23").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 2
$ws.Range("E23").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("C23").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 2
$ws.Range("E23").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100

# ------------------------------------------------------------------
# Cells that flip from a genuine number back to shared-string text:
# copy the text style/format+value from a same-style neighbour that
# already holds the exact text we need (no retyping needed).
# ------------------------------------------------------------------
$ws.Range("C28").Copy($ws.Range("D28"))
$ws.Range("M28").Copy($ws.Range("E28"))

# ------------------------------------------------------------------
# Plain numeric value updates (style/format unchanged).
# ------------------------------------------------------------------
$ws.Range("N14").Value = -50
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -57.894736842105
$ws.Range("I16").Value = 81
$ws.Range("J16").Value = 102
$ws.Range("K16").Value = -20.588235294117
$ws.Range("L16").Value = -39.552238805970
$ws.Range("M16").Value = -54.748603351955
$ws.Range("N16").Value = -87.019230769230
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 22.222222222222
$ws.Range("I17").Value = 196
$ws.Range("J17").Value = 154
$ws.Range("K17").Value = 27.272727272727
$ws.Range("L17").Value = 51.937984496124
$ws.Range("M17").Value = 196.969696969697
$ws.Range("N17").Value = 13.953488372093
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -55.555555555555
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -9.523809523809
$ws.Range("I18").Value = 148
$ws.Range("J18").Value = 191
$ws.Range("K18").Value = -22.513089005235
$ws.Range("L18").Value = 0.680272108843
$ws.Range("M18").Value = -19.565217391304
$ws.Range("N18").Value = -85.700483091787
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -19.298245614035
$ws.Range("I19").Value = 386
$ws.Range("J19").Value = 419
$ws.Range("K19").Value = -7.875894988066
$ws.Range("L19").Value = -15.904139433551
$ws.Range("M19").Value = 24.919093851132
$ws.Range("N19").Value = -7.211538461538
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -14.285714285714
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = -10
$ws.Range("I20").Value = 244
$ws.Range("J20").Value = 255
$ws.Range("K20").Value = -4.313725490196
$ws.Range("L20").Value = 44.378698224852
$ws.Range("M20").Value = 62.666666666666
$ws.Range("N20").Value = -92.690233672858
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -35
$ws.Range("F21").Value = 124
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = -15.068493150684
$ws.Range("I21").Value = 1072
$ws.Range("J21").Value = 1133
$ws.Range("K21").Value = -5.383936451897
$ws.Range("L21").Value = 1.707779886148
$ws.Range("M21").Value = 19.111111111111
$ws.Range("N21").Value = -80.908281389136
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -83.333333333333
$ws.Range("J22").Value = 31
$ws.Range("K22").Value = -19.354838709677
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 40
$ws.Range("J23").Value = 58
$ws.Range("K23").Value = -31.034482758620
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 73.913043478260
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 3.703703703703
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = -22.580645161290
$ws.Range("I24").Value = 824
$ws.Range("J24").Value = 994
$ws.Range("K24").Value = -17.102615694165
$ws.Range("L24").Value = -16.935483870967
$ws.Range("M24").Value = 21.713441654357
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -20
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 59
$ws.Range("H25").Value = -27.118644067796
$ws.Range("I25").Value = 350
$ws.Range("J25").Value = 373
$ws.Range("K25").Value = -6.166219839142
$ws.Range("L25").Value = 16.666666666666
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 9.090909090909
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = 39.393939393939
$ws.Range("I26").Value = 356
$ws.Range("J26").Value = 324
$ws.Range("K26").Value = 9.876543209876
$ws.Range("L26").Value = 13.375796178343
$ws.Range("M26").Value = 24.912280701754
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 24
$ws.Range("J27").Value = 21
$ws.Range("K27").Value = 14.285714285714
$ws.Range("L27").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -80
$ws.Range("L28").Value = -35.714285714285
$ws.Range("L31").Value = 66.666666666666
